$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("A1")
$rng.Font.Color = 255
Write-Host "done"
